# Daily scrape refresh - 2025-12-11 03:35:18 UTC
# Updates Sheet1 rows 2-5 with fresh scrape data, removes the now-unused
# "premium highlight" (yellow fill) from E2, narrows/widens a few columns,
# and drops the two oldest rows (6 & 7) so the sheet shrinks to A1:H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    # Numeric-looking strings (opportunity IDs) would otherwise be
    # auto-coerced to a number by Excel's input parser; a leading
    # apostrophe forces them to stay plain text, matching the scraped
    # source data (which stores ids as text).
    if ($text -match '^[0-9]+$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# --- Row 2 -----------------------------------------------------------
Set-TextValue $ws.Cells.Item(2, 1) "1330536"
Set-TextValue $ws.Cells.Item(2, 2) "https://aiesec.org/opportunity/global-talent/1330536"
Set-TextValue $ws.Cells.Item(2, 3) "[IMPACT FORTALEZA] INSIDE SALES"
Set-TextValue $ws.Cells.Item(2, 4) "Castanhal, PA, Brasil"
Set-TextValue $ws.Cells.Item(2, 5) "No"
$ws.Cells.Item(2, 5).Style = "Normal"
Set-TextValue $ws.Cells.Item(2, 6) "0 applicants"
Set-TextValue $ws.Cells.Item(2, 7) "6 - 18 Months"
Set-TextValue $ws.Cells.Item(2, 8) "Petruz Fruity"
$ws.Cells.Item(2, 1).Style = "Normal"

# --- Row 3 -----------------------------------------------------------
Set-TextValue $ws.Cells.Item(3, 1) "1330505"
Set-TextValue $ws.Cells.Item(3, 2) "https://aiesec.org/opportunity/global-talent/1330505"
Set-TextValue $ws.Cells.Item(3, 3) "Social Media Marketing"
Set-TextValue $ws.Cells.Item(3, 4) "Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye"
Set-TextValue $ws.Cells.Item(3, 5) "No"
Set-TextValue $ws.Cells.Item(3, 6) "7 applicants"
Set-TextValue $ws.Cells.Item(3, 7) "6 - 18 Months"
Set-TextValue $ws.Cells.Item(3, 8) "HATTENA TEKSTİL"
$ws.Cells.Item(3, 1).Style = "Normal"

# --- Row 4 -----------------------------------------------------------
Set-TextValue $ws.Cells.Item(4, 1) "1330498"
Set-TextValue $ws.Cells.Item(4, 2) "https://aiesec.org/opportunity/global-talent/1330498"
Set-TextValue $ws.Cells.Item(4, 3) "Front Office Executive"
Set-TextValue $ws.Cells.Item(4, 4) "Ella, Sri Lanka"
Set-TextValue $ws.Cells.Item(4, 5) "No"
Set-TextValue $ws.Cells.Item(4, 6) "3 applicants"
Set-TextValue $ws.Cells.Item(4, 7) "3 - 6 Months"
Set-TextValue $ws.Cells.Item(4, 8) "Area 4 Eco Cubes - Ella"
$ws.Cells.Item(4, 1).Style = "Normal"

# --- Row 5 -----------------------------------------------------------
Set-TextValue $ws.Cells.Item(5, 1) "1328774"
Set-TextValue $ws.Cells.Item(5, 2) "https://aiesec.org/opportunity/global-talent/1328774"
Set-TextValue $ws.Cells.Item(5, 3) "Digital Marketing"
Set-TextValue $ws.Cells.Item(5, 4) "Leiria, Portugal"
Set-TextValue $ws.Cells.Item(5, 5) "No"
Set-TextValue $ws.Cells.Item(5, 6) "67 applicants"
Set-TextValue $ws.Cells.Item(5, 7) "9 - 12 Weeks"
Set-TextValue $ws.Cells.Item(5, 8) "Multidrive"
$ws.Cells.Item(5, 1).Style = "Normal"

# --- Drop the two stale rows (old rows 6 & 7) -------------------------
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# --- Column width tweaks ----------------------------------------------
# `ColumnWidth` (chars) round-trips into the OOXML `width` attribute with
# a fixed +5/6 padding baked in by this host, so back it out here to land
# on the exact integer widths the sheet should end up with.
$colWidthPad = 5.0 / 6.0
$ws.Columns.Item(3).ColumnWidth = 34 - $colWidthPad
$ws.Columns.Item(4).ColumnWidth = 57 - $colWidthPad
$ws.Columns.Item(6).ColumnWidth = 16 - $colWidthPad
$ws.Columns.Item(8).ColumnWidth = 26 - $colWidthPad
